# Updates cryptocurrency price (D) and 1h volume change (E) figures
# for the rows whose source data refreshed, per the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.258.97"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "3.319.67"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'586.08"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "'183.67"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "'0.648"
$ws.Range("E7").Value = "  +6.97%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "3.899.94"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").Value = "66.294.49"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'26.13"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "3.323.47"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "'424.33"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D20").Value = "'13.19"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "'71.81"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'5.68"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "3.464.09"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'0.514"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "'0.202"
$ws.Range("E27").Value = "  +5.91%  "
$ws.Range("E28").Value = "  -3.20%  "
$ws.Range("D29").Value = "'8.92"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'1.91"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'22.40"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D34").Value = "'5.18"
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("D35").Value = "'6.58"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("D37").Value = "'159.92"
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").Value = "2.889.53"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "'26.47"
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("D42").Value = "'0.764"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").Value = "'40.04"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "'0.0664"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'5.94"
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("D47").Value = "'2.30"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").Value = "'23.24"
$ws.Range("E48").Value = "  -5.03%  "
$ws.Range("D49").Value = "'313.76"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  +5.03%  "
